$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.530.62"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "3.593.39"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.49"
$ws.Range("E5").Value = "  +7.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "557.24"
$ws.Range("E6").Value = "  -5.12%  "
$ws.Range("D7").Value = "3.592.21"
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.613"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.07"
$ws.Range("E11").Value = "  +9.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.152"
$ws.Range("E12").Value = "  +3.95%  "
$ws.Range("E13").Value = "  +11.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.97"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "4.182.31"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "3.596.56"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.92"
$ws.Range("E18").Value = "  +3.50%  "
$ws.Range("D19").Value = "67.506.44"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.24"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.07"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "398.74"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.88"
$ws.Range("E23").Value = "  +16.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.12"
$ws.Range("E24").Value = "  -5.21%  "
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.47"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.85"
$ws.Range("E28").Value = "  +8.74%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.11"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.43"
$ws.Range("E30").Value = "  +19.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.43"
$ws.Range("E31").Value = "  +5.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.40"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "664.73"
$ws.Range("E33").Value = "  +7.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.16"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.60"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.25"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.430"
$ws.Range("E38").Value = "  +11.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "0.0₃0770"
$ws.Range("E40").Value = "  +2.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.20"
$ws.Range("E41").Value = "  +14.28%  "
$ws.Range("D42").Value = "3.250.80"
$ws.Range("E42").Value = "  +8.80%  "
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("E44").Value = "  +13.51%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.01"
$ws.Range("E45").Value = "  +29.41%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.73"
$ws.Range("E48").Value = "  +10.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.14"
$ws.Range("E49").Value = "  +4.35%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.72"
$ws.Range("E51").Value = "  +1.57%  "
